$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates: issue number and week-covering dates
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

function Set-TextCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.NumberFormat = "General"
}

function Set-CountCell($addr, $num) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "#,##0"
    $r.Value = $num
}

function Set-PctCell($addr, $num) {
    $r = $ws.Range($addr)
    $r.NumberFormat = '#,##0.0;"-"#,##0.0'
    $r.Value = $num
}

# Row 15
Set-TextCell "C15" "0"
Set-TextCell "D15" "0"
Set-TextCell "E15" "***.*"

# Row 16
Set-CountCell "C16" 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("I16").Value = 52
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 4
$ws.Range("L16").Value = -24.637681159420
$ws.Range("M16").Value = -48
$ws.Range("N16").Value = -90.316573556797

# Row 17
Set-CountCell "D17" 1
Set-PctCell "E17" 0
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 96
$ws.Range("J17").Value = 86
$ws.Range("K17").Value = 11.627906976744
$ws.Range("L17").Value = -25.581395348837
$ws.Range("M17").Value = -4
$ws.Range("N17").Value = -65.467625899280

# Row 18
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 91
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = -24.166666666666
$ws.Range("L18").Value = -33.088235294117
$ws.Range("M18").Value = -61.764705882352
$ws.Range("N18").Value = -91.853178155774

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 11.764705882352
$ws.Range("I19").Value = 488
$ws.Range("J19").Value = 362
$ws.Range("K19").Value = 34.806629834254
$ws.Range("L19").Value = 53.459119496855
$ws.Range("M19").Value = 55.910543130990
$ws.Range("N19").Value = 5.856832971800

# Row 20
Set-TextCell "D20" "0"
Set-TextCell "E20" "***.*"
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -11.111111111111
$ws.Range("I20").Value = 102
$ws.Range("K20").Value = 18.604651162790
$ws.Range("L20").Value = -23.880597014925
$ws.Range("M20").Value = -29.655172413793
$ws.Range("N20").Value = -94.606028556319

# Row 21
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 23.076923076923
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 60
$ws.Range("H21").Value = 5
$ws.Range("I21").Value = 842
$ws.Range("J21").Value = 720
$ws.Range("K21").Value = 16.944444444444
$ws.Range("L21").Value = 4.596273291925
$ws.Range("M21").Value = -6.755260243632
$ws.Range("N21").Value = -80.491195551436

# Row 22
Set-CountCell "D22" 1
Set-PctCell "E22" -100
Set-TextCell "F22" "0"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = 57.142857142857

# Row 24
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 153
$ws.Range("G24").Value = 123
$ws.Range("H24").Value = 24.390243902439
$ws.Range("I24").Value = 1826
$ws.Range("J24").Value = 1081
$ws.Range("K24").Value = 68.917668825161
$ws.Range("L24").Value = 87.860082304526
$ws.Range("M24").Value = 55.140186915887

# Row 25
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = 20
$ws.Range("I25").Value = 313
$ws.Range("J25").Value = 245
$ws.Range("K25").Value = 27.755102040816
$ws.Range("L25").Value = 37.885462555066
$ws.Range("M25").Value = -17.847769028871

# Row 26
Set-TextCell "C26" "0"
Set-TextCell "D26" "0"
Set-TextCell "E26" "***.*"

# Row 27
Set-TextCell "C27" "0"
Set-TextCell "D27" "0"
Set-TextCell "E27" "***.*"
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -50

# Row 30
Set-TextCell "D30" "0"
Set-TextCell "E30" "***.*"
